$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 294.66666
$ws.Range("I12").Value = 294.66666
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 294.66666
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -124.66666
$ws.Range("N12").ClearContents()
$ws.Range("H18").Value = 1239.7778
$ws.Range("I18").Value = 1481.7142
$ws.Range("J18").Value = 1085.8182
$ws.Range("K18").Value = 1481.7142
$ws.Range("L18").Value = 1085.8182
$ws.Range("M18").Value = -1197.7142
$ws.Range("N18").Value = -1653.8182
$ws.Range("H28").Value = 1629.875
$ws.Range("I28").Value = 1907
$ws.Range("J28").Value = 798.5
$ws.Range("K28").Value = 1907
$ws.Range("L28").Value = 798.5
$ws.Range("M28").Value = -1422
$ws.Range("N28").Value = -1768.5
$ws.Range("H40").Value = 1870.7368
$ws.Range("I40").Value = 1622.037
$ws.Range("J40").Value = 2481.182
$ws.Range("K40").Value = 1622.037
$ws.Range("L40").Value = 2481.182
$ws.Range("M40").Value = -1447.037
$ws.Range("N40").Value = -2831.182
$ws.Range("H41").Value = 532
$ws.Range("I41").Value = 546.4
$ws.Range("J41").Value = 524
$ws.Range("K41").Value = 546.4
$ws.Range("L41").Value = 524
$ws.Range("M41").Value = -106.4
$ws.Range("N41").Value = -1404
$ws.Range("H53").Value = 354.66666
$ws.Range("I53").Value = 257.8889
$ws.Range("J53").Value = 499.83334
$ws.Range("K53").Value = 257.8889
$ws.Range("L53").Value = 499.83334
$ws.Range("M53").Value = 379.1111
$ws.Range("N53").Value = -1773.83334
$ws.Range("H96").Value = 13057.889
$ws.Range("I96").Value = 26756.5
$ws.Range("J96").Value = 2099
$ws.Range("K96").Value = 80269.5
$ws.Range("L96").Value = 6297
$ws.Range("M96").Value = -78896.5
$ws.Range("N96").Value = -9043
$ws.Range("H98").Value = 1749.75
$ws.Range("I98").Value = 1749.75
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 1749.75
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -251.75
$ws.Range("H108").Value = 0
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
$ws.Range("H122").Value = 1749.75
$ws.Range("I122").Value = 1749.75
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5249.25
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2799.25
$ws.Range("H129").Value = 4542.273
$ws.Range("I129").Value = 9999
$ws.Range("J129").Value = 3329.6667
$ws.Range("K129").Value = 29997
$ws.Range("L129").Value = 9989.000100000001
$ws.Range("M129").Value = -24997
$ws.Range("N129").Value = -19989.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3350852.8
$ws.Range("I32").Value = 3336581.5
$ws.Range("J32").Value = 3500700
$ws.Range("K32").Value = 3336581.5
$ws.Range("L32").Value = 3500700
$ws.Range("M32").Value = -3336294.5
$ws.Range("N32").Value = -3501274
$ws.Range("H46").Value = 4805.5713
$ws.Range("I46").Value = 4767
$ws.Range("J46").Value = 4812
$ws.Range("K46").Value = 4767
$ws.Range("L46").Value = 4812
$ws.Range("M46").Value = -4448
$ws.Range("N46").Value = -5450
$ws.Range("H61").Value = 2718.1667
$ws.Range("I61").Value = 2261.8
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 2261.8
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -2049.8
$ws.Range("N61").Value = -5424
$ws.Range("H88").Value = 2575.8333
$ws.Range("I88").Value = 3506
$ws.Range("J88").Value = 2389.8
$ws.Range("K88").Value = 3506
$ws.Range("L88").Value = 2389.8
$ws.Range("M88").Value = -3100
$ws.Range("N88").Value = -3201.8
$ws.Range("H91").Value = 2575.8333
$ws.Range("I91").Value = 3506
$ws.Range("J91").Value = 2389.8
$ws.Range("K91").Value = 3506
$ws.Range("L91").Value = 2389.8
$ws.Range("M91").Value = -2102
$ws.Range("N91").Value = -5197.8
$ws.Range("H97").Value = 1015.2857
$ws.Range("I97").Value = 937
$ws.Range("J97").Value = 1485
$ws.Range("K97").Value = 937
$ws.Range("L97").Value = 1485
$ws.Range("M97").Value = -441
$ws.Range("N97").Value = -2477
$ws.Range("H110").Value = 1143.7646
$ws.Range("I110").Value = 1082.5
$ws.Range("J110").Value = 2124
$ws.Range("K110").Value = 1082.5
$ws.Range("L110").Value = 2124
$ws.Range("M110").Value = 962.5
$ws.Range("N110").Value = -6214
$ws.Range("H136").Value = 2718.1667
$ws.Range("I136").Value = 2261.8
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 6785.400000000001
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -4235.400000000001
$ws.Range("N136").Value = -20100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 7750
$ws.Range("I20").Value = 5500
$ws.Range("J20").Value = 10000
$ws.Range("K20").Value = 5500
$ws.Range("L20").Value = 10000
$ws.Range("M20").Value = -5253
$ws.Range("N20").Value = -10494
$ws.Range("H107").Value = 1427.8
$ws.Range("I107").Value = 1380.8334
$ws.Range("J107").Value = 1498.25
$ws.Range("K107").Value = 1380.8334
$ws.Range("L107").Value = 1498.25
$ws.Range("M107").Value = 539.1666
$ws.Range("N107").Value = -5338.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2222
$ws.Range("I31").Value = 2222
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 2222
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -1927
$ws.Range("H34").Value = 2222
$ws.Range("I34").Value = 2222
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 2222
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -2020
$ws.Range("H86").Value = 19997
$ws.Range("I86").Value = 19997
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 19997
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -18874
$ws.Range("H89").Value = 19997
$ws.Range("I89").Value = 19997
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 99985
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -94369
$ws.Range("H94").Value = 114655.4
$ws.Range("I94").Value = 223732.4
$ws.Range("J94").Value = 5578.4
$ws.Range("K94").Value = 223732.4
$ws.Range("L94").Value = 5578.4
$ws.Range("M94").Value = -223281.4
$ws.Range("N94").Value = -6480.4
$ws.Range("H132").Value = 4475
$ws.Range("I132").Value = 4594.5
$ws.Range("J132").Value = 3997
$ws.Range("K132").Value = 13783.5
$ws.Range("L132").Value = 11991
$ws.Range("M132").Value = -11253.5
$ws.Range("N132").Value = -17051
$ws.Range("H141").Value = 40113
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 40113
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 40113
$ws.Range("N141").Value = -50473

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 13588899
$ws.Range("I4").Value = 15187511
$ws.Range("J4").Value = 800000
$ws.Range("K4").Value = 45562533
$ws.Range("L4").Value = 2400000
$ws.Range("M4").Value = -45562421
$ws.Range("N4").Value = -2400224
$ws.Range("H95").Value = 8893.5
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 8893.5
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 26680.5
$ws.Range("N95").Value = -30798.5
$ws.Range("H113").Value = 1648.5
$ws.Range("I113").Value = 1647.5
$ws.Range("J113").Value = 1649.5
$ws.Range("K113").Value = 4942.5
$ws.Range("L113").Value = 4948.5
$ws.Range("M113").Value = -2772.5
$ws.Range("N113").Value = -9288.5
$ws.Range("H122").Value = 741.5
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 741.5
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 6673.5
$ws.Range("N122").Value = -11573.5
$ws.Range("H131").Value = 716239.9
$ws.Range("I131").Value = 1319
$ws.Range("J131").Value = 1252430.5
$ws.Range("K131").Value = 3957
$ws.Range("L131").Value = 3757291.5
$ws.Range("M131").Value = 1083
$ws.Range("N131").Value = -3767371.5
$ws.Range("H137").Value = 2120.2727
$ws.Range("I137").Value = 733
$ws.Range("J137").Value = 2640.5
$ws.Range("K137").Value = 2199
$ws.Range("L137").Value = 7921.5
$ws.Range("M137").Value = 2901
$ws.Range("N137").Value = -18121.5
$ws.Range("H138").Value = 2171.75
$ws.Range("I138").Value = 2096
$ws.Range("J138").Value = 2399
$ws.Range("K138").Value = 6288
$ws.Range("L138").Value = 7197
$ws.Range("M138").Value = -1148
$ws.Range("N138").Value = -17477

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 161687.86
$ws.Range("I10").Value = 187502.5
$ws.Range("J10").Value = 6800
$ws.Range("K10").Value = 187502.5
$ws.Range("L10").Value = 6800
$ws.Range("M10").Value = -187333.5
$ws.Range("N10").Value = -7138
$ws.Range("H132").Value = 6493.963
$ws.Range("I132").Value = 6680.9585
$ws.Range("J132").Value = 4998
$ws.Range("K132").Value = 20042.8755
$ws.Range("L132").Value = 14994
$ws.Range("M132").Value = -17512.8755
$ws.Range("N132").Value = -20054
$ws.Range("H134").Value = 24398.8
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 24398.8
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 73196.39999999999
$ws.Range("N134").Value = -78266.39999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3521.889
$ws.Range("I22").Value = 4416.6665
$ws.Range("J22").Value = 1732.3334
$ws.Range("K22").Value = 4416.6665
$ws.Range("L22").Value = 1732.3334
$ws.Range("M22").Value = -4121.6665
$ws.Range("N22").Value = -2322.3334
$ws.Range("H27").Value = 3521.889
$ws.Range("I27").Value = 4416.6665
$ws.Range("J27").Value = 1732.3334
$ws.Range("K27").Value = 4416.6665
$ws.Range("L27").Value = 1732.3334
$ws.Range("M27").Value = -4309.6665
$ws.Range("N27").Value = -1946.3334
$ws.Range("H61").Value = 2799
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 2799
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 2799
$ws.Range("N61").Value = -3203
$ws.Range("H108").Value = 50000
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 50000
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 50000
$ws.Range("N108").Value = -57680
$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
$ws.Range("H113").Value = 2799
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 2799
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 2799
$ws.Range("N113").Value = -7139
$ws.Range("H122").Value = 6320.3213
$ws.Range("I122").Value = 4826.273
$ws.Range("J122").Value = 7287.0586
$ws.Range("K122").Value = 14478.819
$ws.Range("L122").Value = 21861.1758
$ws.Range("M122").Value = -12028.819
$ws.Range("N122").Value = -26761.1758
$ws.Range("H132").Value = 3565.3333
$ws.Range("I132").Value = 3632.6667
$ws.Range("J132").Value = 3498
$ws.Range("K132").Value = 10898.0001
$ws.Range("L132").Value = 10494
$ws.Range("M132").Value = -8368.000100000001
$ws.Range("N132").Value = -15554
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 47297.5
$ws.Range("I70").Value = 49595
$ws.Range("J70").Value = 45000
$ws.Range("K70").Value = 49595
$ws.Range("L70").Value = 45000
$ws.Range("M70").Value = -49280
$ws.Range("N70").Value = -45630
$ws.Range("H73").Value = 47297.5
$ws.Range("I73").Value = 49595
$ws.Range("J73").Value = 45000
$ws.Range("K73").Value = 49595
$ws.Range("L73").Value = 45000
$ws.Range("M73").Value = -48503
$ws.Range("N73").Value = -47184
$ws.Range("H132").Value = 3238.8
$ws.Range("I132").Value = 3238.8
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 9716.400000000001
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -7186.400000000001
$ws.Range("H136").Value = 1981.2
$ws.Range("I136").Value = 1226.5
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 3679.5
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -1129.5
$ws.Range("N136").Value = -20100
